$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.8908616666666668
$ws.Range("N2").Value = 2.672585
$ws.Range("O2").Value = 0.04079002072021364
$ws.Range("P2").Value = 0.04079002072021363
$ws.Range("Q2").Value = 9.29546740633889
$ws.Range("R2").Value = 83.65920665704999
$ws.Range("S2").Value = 0.03961156430308091
$ws.Range("T2").Value = 0.03961156430308091

$ws.Range("M3").Value = 16.81477433333333
$ws.Range("O3").Value = 0.7699006693471485
$ws.Range("P3").Value = 0.7699006693471484
$ws.Range("S3").Value = 0.7476576214563364
$ws.Range("T3").Value = 0.7476576214563364

$ws.Range("M4").Value = 3.879966
$ws.Range("N4").Value = 11.639898
$ws.Range("O4").Value = 0.1776526024808091
$ws.Range("P4").Value = 0.1776526024808091
$ws.Range("Q4").Value = 40.48450936906
$ws.Range("R4").Value = 364.36058432154
$ws.Range("S4").Value = 0.1725200762962835
$ws.Range("T4").Value = 0.1725200762962835

$ws.Range("M5").Value = 0.2545846666666667
$ws.Range("N5").Value = 0.763754
$ws.Range("O5").Value = 0.01165670745182886
$ws.Range("P5").Value = 0.01165670745182886
$ws.Range("Q5").Value = 2.656398360935555
$ws.Range("R5").Value = 23.90758524842
$ws.Range("S5").Value = 0.01131993582345753
$ws.Range("T5").Value = 0.01131993582345753

$ws.Range("M6").Value = 0.8908616666666668
$ws.Range("N6").Value = 2.672585
$ws.Range("O6").Value = 0.04079002072021364
$ws.Range("P6").Value = 0.04079002072021363
$ws.Range("Q6").Value = 0.27654306029
$ws.Range("R6").Value = 2.48888754261
$ws.Range("S6").Value = 0.001178456417132721
$ws.Range("T6").Value = 0.001178456417132721

$ws.Range("M7").Value = 16.81477433333333
$ws.Range("O7").Value = 0.7699006693471485
$ws.Range("P7").Value = 0.7699006693471484
$ws.Range("Q7").Value = 5.219675878102
$ws.Range("R7").Value = 46.97708290291799
$ws.Range("S7").Value = 0.02224304789081197
$ws.Range("T7").Value = 0.02224304789081197

$ws.Range("M8").Value = 3.879966
$ws.Range("N8").Value = 11.639898
$ws.Range("O8").Value = 0.1776526024808091
$ws.Range("P8").Value = 0.1776526024808091
$ws.Range("Q8").Value = 1.204426805652
$ws.Range("R8").Value = 10.839841250868
$ws.Range("S8").Value = 0.005132526184525591
$ws.Range("T8").Value = 0.00513252618452559

$ws.Range("M9").Value = 0.2545846666666667
$ws.Range("N9").Value = 0.763754
$ws.Range("O9").Value = 0.01165670745182886
$ws.Range("P9").Value = 0.01165670745182886
$ws.Range("Q9").Value = 0.079028681396
$ws.Range("R9").Value = 0.7112581325639999
$ws.Range("S9").Value = 0.0003367716283713275
$ws.Range("T9").Value = 0.0003367716283713275
